$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recalculated values in existing rows (row 2) ---
$ws.Range("M2").Value = 0.8639135
$ws.Range("N2").Value = 1.727827
$ws.Range("O2").Value = 0.1895490737713731
$ws.Range("P2").Value = 0.1380033232738433
$ws.Range("Q2").Value = 0.10473525730225
$ws.Range("R2").Value = 0.418941029209
$ws.Range("S2").Value = 0.1895490737713731
$ws.Range("T2").Value = 0.1380033232738433

# --- Row 3 ---
$ws.Range("O3").Value = 0.7425593442349591
$ws.Range("P3").Value = 0.8109429541930055
$ws.Range("Q3").Value = 0.4103008388974999
$ws.Range("R3").Value = 2.461805033384999
$ws.Range("S3").Value = 0.7425593442349591
$ws.Range("T3").Value = 0.8109429541930055

# --- Row 4 ---
$ws.Range("M4").Value = 0.2890925
$ws.Range("N4").Value = 0.5781849999999999
$ws.Range("O4").Value = 0.06342905349812297
$ws.Range("P4").Value = 0.04618023185601746
$ws.Range("Q4").Value = 0.03504769559875
$ws.Range("R4").Value = 0.140190782395
$ws.Range("S4").Value = 0.06342905349812297
$ws.Range("T4").Value = 0.04618023185601746

# --- Row 5 ---
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01599833333333333
$ws.Range("N5").Value = 0.047995
$ws.Range("O5").Value = 0.003510153811948785
$ws.Range("P5").Value = 0.003833410116017465
$ws.Range("Q5").Value = 0.001939533944166667
$ws.Range("R5").Value = 0.011637203665
$ws.Range("S5").Value = 0.003510153811948785
$ws.Range("T5").Value = 0.003833410116017465

# --- New row 6 ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nrtn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.1212335
$ws.Range("H6").Value = 0.242467
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.004340666666666667
$ws.Range("N6").Value = 0.013022
$ws.Range("O6").Value = 0.0009523746835961471
$ws.Range("P6").Value = 0.001040080561116354
$ws.Range("Q6").Value = 0.0005262342123333334
$ws.Range("R6").Value = 0.003157405274
$ws.Range("S6").Value = 0.0009523746835961471
$ws.Range("T6").Value = 0.001040080561116354
